$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header date strings (keep as text, preserve original cell style) ---
function Set-TextValue($cell, $text) {
    $fmt = $cell.NumberFormat
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.NumberFormat = $fmt
}

Set-TextValue $ws.Range("A5") "2024/04/01"
Set-TextValue $ws.Range("C5") "2025/04/01"

# --- Numeric value updates (columns A-D across affected rows) ---
$ws.Range("A9").Value = 88.177612248629
$ws.Range("B9").Value = 593.311
$ws.Range("A10").Value = 90.530797225846
$ws.Range("B10").Value = 26.368
$ws.Range("C10").Value = 49.385428826478
$ws.Range("D10").Value = 14.384
$ws.Range("A11").Value = 82.48569612206
$ws.Range("B11").Value = 18.165
$ws.Range("C11").Value = 71.478521478521
$ws.Range("D11").Value = 15.741
$ws.Range("A12").Value = 25.021085933932
$ws.Range("B12").Value = 75.054
$ws.Range("A14").Value = 87.671690072142
$ws.Range("B14").Value = 34.149
$ws.Range("A15").Value = 66.969774590164
$ws.Range("B15").Value = 15.687
$ws.Range("C15").Value = 51.263661202186
$ws.Range("D15").Value = 12.008
$ws.Range("A16").Value = 24.705283690951
$ws.Range("B16").Value = 2.913
$ws.Range("C16").Value = 35.009753201594
$ws.Range("D16").Value = 4.128
$ws.Range("A17").Value = 15.412950756107
$ws.Range("B17").Value = 0.795
$ws.Range("C17").Value = 24.893369523071
$ws.Range("D17").Value = 1.284
$ws.Range("A19").Value = 28.620930972943
$ws.Range("B19").Value = 137.4564
$ws.Range("C19").Value = 23.810631332307
$ws.Range("D19").Value = 114.3542
$ws.Range("C21").Value = 55.340091589675
$ws.Range("D21").Value = 104.444688
$ws.Range("A22").Value = 12.895200369104
$ws.Range("B22").Value = 50.588
$ws.Range("A23").Value = 15.948002536462
$ws.Range("B23").Value = 2.012
$ws.Range("C23").Value = 30.976537729867
$ws.Range("D23").Value = 3.908
$ws.Range("A24").Value = 15.424927476987
$ws.Range("B24").Value = 36.86557667
$ws.Range("C24").Value = 82.36041634481
$ws.Range("D24").Value = 135.4418708
$ws.Range("A25").Value = 58.540218470705
$ws.Range("B25").Value = 4.716
$ws.Range("C25").Value = 60.329569519546
$ws.Range("D25").Value = 2.820595
$ws.Range("A26").Value = 15.569302321896
$ws.Range("B26").Value = 7.8601
$ws.Range("A28").Value = 23.958932754466
$ws.Range("B28").Value = 270.631
$ws.Range("C28").Value = 35.634166163522
$ws.Range("D28").Value = 402.51
$ws.Range("A29").Value = 94.860450887184
$ws.Range("B29").Value = 60.465
$ws.Range("C29").Value = 97.202742347939
$ws.Range("D29").Value = 61.958
$ws.Range("A30").Value = 41.349234603062
$ws.Range("B30").Value = 27.876
$ws.Range("C30").Value = 42.47063011748
$ws.Range("D30").Value = 28.632
$ws.Range("A31").Value = 91.930972478897
$ws.Range("B31").Value = 36.811
$ws.Range("C31").Value = 96.908246341342
$ws.Range("D31").Value = 38.804
$ws.Range("A32").Value = 31.630649237835
$ws.Range("B32").Value = 68.456
$ws.Range("A33").Value = 53.234067295582
$ws.Range("B33").Value = 1875.038
$ws.Range("C33").Value = 57.106249070197
$ws.Range("D33").Value = 2011.426
$ws.Range("A34").Value = 50.041550793391
$ws.Range("B34").Value = 76.476
$ws.Range("A36").Value = 41.446235662074
$ws.Range("B36").Value = 116.747
$ws.Range("A37").Value = 52.700824882476
$ws.Range("B37").Value = 17.825
$ws.Range("C37").Value = 66.620347101085
$ws.Range("D37").Value = 22.533
$ws.Range("A39").Value = 27.374977943922
$ws.Range("B39").Value = 266.848
$ws.Range("C39").Value = 61.826571521192
$ws.Range("D39").Value = 602.678
$ws.Range("A40").Value = 31.143140886276
$ws.Range("B40").Value = 11.603
$ws.Range("A41").Value = 35.586871209418
$ws.Range("B41").Value = 19.95
$ws.Range("C41").Value = 41.787370674278
$ws.Range("D41").Value = 23.426
$ws.Range("A42").Value = 11.592655849701
$ws.Range("B42").Value = 1.629
$ws.Range("A43").Value = 4.8540187497755
$ws.Range("B43").Value = 32.433
$ws.Range("C43").Value = 18.623160642234
$ws.Range("D43").Value = 124.434
$ws.Range("A44").Value = 1.9386577232232
$ws.Range("B44").Value = 51.51
$ws.Range("C44").Value = 4.3797631382544
$ws.Range("D44").Value = 116.37
$ws.Range("A45").Value = 24.611715062829
$ws.Range("B45").Value = 58.093
$ws.Range("C45").Value = 20.815716113507
$ws.Range("D45").Value = 49.133
$ws.Range("C46").Value = 87.423438138016
$ws.Range("D46").Value = 2.141
$ws.Range("A47").Value = 48.789915966387
$ws.Range("B47").Value = 69.672
$ws.Range("C47").Value = 54.572128851541
$ws.Range("D47").Value = 77.929
$ws.Range("A48").Value = 68.929889298893
$ws.Range("B48").Value = 0.934
$ws.Range("C48").Value = 89.815498154982
$ws.Range("D48").Value = 1.217
$ws.Range("A49").Value = 19.080502908613
$ws.Range("B49").Value = 10.168
$ws.Range("C49").Value = 17.461062112967
$ws.Range("D49").Value = 9.305
$ws.Range("A50").Value = 72.411558895944
$ws.Range("B50").Value = 11.176
$ws.Range("C50").Value = 94.674096151354
$ws.Range("D50").Value = 14.612
$ws.Range("A51").Value = 86.234695201973
$ws.Range("B51").Value = 59.796
$ws.Range("C51").Value = 60.648101411863
$ws.Range("D51").Value = 42.054
$ws.Range("A52").Value = 44.325574895579
$ws.Range("B52").Value = 28.547
$ws.Range("C52").Value = 65.833889725634
$ws.Range("D52").Value = 42.399
$ws.Range("A53").Value = 25.585692728714
$ws.Range("B53").Value = 80.03
$ws.Range("C53").Value = 64.901276247474
$ws.Range("D53").Value = 203.006
$ws.Range("A54").Value = 14.682426563177
$ws.Range("B54").Value = 65.413
$ws.Range("C54").Value = 43.768728157497
$ws.Range("D54").Value = 194.998
$ws.Range("A56").Value = 90.789473684211
$ws.Range("B56").Value = 0.207
$ws.Range("C56").Value = 61.842105263158
$ws.Range("D56").Value = 0.141
$ws.Range("A57").Value = 26.193892191023
$ws.Range("B57").Value = 2.556
$ws.Range("A58").Value = 16.121242610992
$ws.Range("B58").Value = 6.409
$ws.Range("A59").Value = 22.623757924554
$ws.Range("B59").Value = 20.127
$ws.Range("C59").Value = 63.924733600108
$ws.Range("D59").Value = 56.87
$ws.Range("A60").Value = 12.330188047915
$ws.Range("B60").Value = 36.85
$ws.Range("A61").Value = 27.892929315238
$ws.Range("B61").Value = 25.259
$ws.Range("C61").Value = 50.190487759091
$ws.Range("D61").Value = 45.451
$ws.Range("A62").Value = 38.184159690921
$ws.Range("B62").Value = 1.779
$ws.Range("C62").Value = 21.764327108822
$ws.Range("D62").Value = 1.014
$ws.Range("C63").Value = 61.909555832964
$ws.Range("D63").Value = 15.374
$ws.Range("A64").Value = 95.70072585148
$ws.Range("B64").Value = 13.712
$ws.Range("C64").Value = 78.601340033501
$ws.Range("D64").Value = 11.262
$ws.Range("A65").Value = 36.902485659656
$ws.Range("B65").Value = 3.474
$ws.Range("A66").Value = 46.019532406037
$ws.Range("B66").Value = 3.11
$ws.Range("C66").Value = 85.557857354247
$ws.Range("D66").Value = 5.782
$ws.Range("A68").Value = 66
$ws.Range("B68").Value = 0.726
$ws.Range("A69").Value = 8.2817497926786
$ws.Range("B69").Value = 100.665
$ws.Range("C69").Value = 11.458868090932
$ws.Range("D69").Value = 139.283
$ws.Range("C70").Value = 35.072882139751
$ws.Range("D70").Value = 78.625073
$ws.Range("A73").Value = 31.607637888777
$ws.Range("B73").Value = 88.5573
$ws.Range("A74").Value = 30.655276651713
$ws.Range("B74").Value = 4942.28977667
$ws.Range("C74").Value = 38.306526422121
$ws.Range("D74").Value = 6421.0301668
